$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new rows at row 2, pushing existing data (old rows 2..176) down to rows 4..178.
# Pull formatting from the row below (the data rows) rather than the header row above.
$ws.Rows("2:3").Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown, [Microsoft.Office.Interop.Excel.XlInsertFormatOrigin]::xlFormatFromRightOrBelow)

# Fill in the two new rows with their transaction data.
$ws.Range("E2").Value = "Deposit"
$ws.Range("N2").Value = "Wiretransfer"
$ws.Range("P2").Value = "Bintense"
$ws.Range("T2").Value = 108567

$ws.Range("E3").Value = "Deposit"
$ws.Range("N3").Value = "Crypto"
$ws.Range("P3").Value = "BTC"
$ws.Range("T3").Value = 26516.7299

# Update the view: scroll so column G is the left-most visible column and
# select T2:T3 (the newly added amounts) as the active selection.
$ws.Application.ActiveWindow.ScrollColumn = 7
$ws.Range("T2:T3").Select()
